$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 2021 row (row 5) with the same style as the preceding year
# label cells (A2:A4) and plain numeric values for the rest of the row.
$ws.Range("A5").Value = "2021年"
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B5").Value = 47.17
$ws.Range("C5").Value = 40.869
$ws.Range("D5").Value = 32.976
$ws.Range("E5").Value = 35.953
$ws.Range("F5").Value = 35.698
$ws.Range("G5").Value = 43.348
$ws.Range("H5").Value = 37.377
$ws.Range("I5").Value = 39.993
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 30.626
$ws.Range("L5").Value = 63.652
$ws.Range("M5").Value = 23.292
